# Adding new TestCase in Notifications
# Appends two new test-case rows (TestCase_E33 / TestCase_E34) to the
# "Test Cases" sheet, right after the existing last row (33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the last existing data row (row 32, which already
# carries the "wrap text" style on column C) onto the two new rows.
$ws.Range("A32:E32").Copy($ws.Range("A34:E34"))
$ws.Range("A32:E32").Copy($ws.Range("A35:E35"))

# Row 34 - write Description/Jira id/TCID first (in that order) so new
# shared-string entries are appended in the same order as the source edit,
# then the Runmode/Results columns which reuse existing shared strings.
$ws.Range("C34").Value = "Verify that anyone can see the public watchlists of a user on user's profile page"
$ws.Range("B34").Value = "OPQA-321"
$ws.Range("A34").Value = "TestCase_E33"
$ws.Range("D34").Value = "Y"
$ws.Range("E34").Value = "PASS"

# Row 35
$ws.Range("C35").Value = "Verify that no one can see the private watchlists of a user on user's profile page"
$ws.Range("B35").Value = "OPQA-329"
$ws.Range("A35").Value = "TestCase_E34"
$ws.Range("D35").Value = "Y"
$ws.Range("E35").Value = "PASS"

# Scroll the view back to the top (column C) and select the full new
# Results range, matching the saved selection state.
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$null = $ws.Range("E2:E35").Select()
